$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "Bi2Ti2O7"
$ws.Cells.Item(2, 2).Value = 3.51
$ws.Cells.Item(3, 1).Value = "Bi1.5Al0.5Ti2O7"
$ws.Cells.Item(3, 2).Value = 3.7
$ws.Cells.Item(4, 1).Value = "Bi1.5Sc0.5Ti2O7"
$ws.Cells.Item(4, 2).Value = 3.73
$ws.Cells.Item(5, 1).Value = "Bi1.5Ga0.5Ti2O7"
$ws.Cells.Item(5, 2).Value = 3.67
$ws.Cells.Item(6, 1).Value = "Bi1.5Y0.5Ti2O7"
$ws.Cells.Item(6, 2).Value = 3.75
$ws.Cells.Item(7, 1).Value = "Bi1.5In0.5Ti2O7"
$ws.Cells.Item(7, 2).Value = 3.54
$ws.Cells.Item(8, 1).Value = "Bi1.5La0.5Ti2O7"
$ws.Cells.Item(8, 2).Value = 3.61
$ws.Cells.Item(9, 1).Value = "Bi4Ti3O12"
$ws.Cells.Item(9, 2).Value = 3.45
$ws.Cells.Item(10, 1).Value = "Bi3.5Al0.5Ti3O12"
$ws.Cells.Item(10, 2).Value = 3.39
$ws.Cells.Item(11, 1).Value = "Bi3.5Ga0.5Ti3O12"
$ws.Cells.Item(11, 2).Value = 3.43
$ws.Cells.Item(12, 1).Value = "Bi3.5In0.5Ti3O12"
$ws.Cells.Item(12, 2).Value = 3.44
$ws.Cells.Item(13, 1).Value = "Bi3.5Sc0.5Ti3O12"
$ws.Cells.Item(13, 2).Value = 3.44
$ws.Cells.Item(14, 1).Value = "Bi3.5Y0.5Ti3O12"
$ws.Cells.Item(14, 2).Value = 3.42
$ws.Cells.Item(15, 1).Value = "Bi3.5La0.5Ti3O12"
$ws.Cells.Item(15, 2).Value = 3.33
$ws.Cells.Item(16, 1).Value = "Bi4Ti2.5Al0.5O12"
$ws.Cells.Item(16, 2).Value = 3.2
$ws.Cells.Item(17, 1).Value = "Bi4Ti2.5Ga0.5O12"
$ws.Cells.Item(17, 2).Value = 3.3
$ws.Cells.Item(18, 1).Value = "Bi4Ti2.5Sc0.5O12"
$ws.Cells.Item(18, 2).Value = 3.33
$ws.Cells.Item(19, 1).Value = "Bi12TiO20"
$ws.Cells.Item(19, 2).Value = 3.07
$ws.Cells.Item(20, 1).Value = "Bi1.5Cr0.5Ti2O7"
$ws.Cells.Item(20, 2).Value = 3.18
$ws.Cells.Item(21, 1).Value = "Bi1.5Fe0.5Ti2O7"
$ws.Cells.Item(21, 2).Value = 0
$ws.Cells.Item(22, 1).Value = "Bi4Ti2O11"
$ws.Cells.Item(22, 2).Value = 3.32
$ws.Cells.Item(23, 1).Value = "Bi1.5Li0.5Ti2O7"
$ws.Cells.Item(23, 2).Value = 3.85
$ws.Cells.Item(24, 1).Value = "Bi1.5Na0.5Ti2O7"
$ws.Cells.Item(24, 2).Value = 3.81
$ws.Cells.Item(25, 1).Value = "Bi1.5Cu0.5Ti2O7"
$ws.Cells.Item(25, 2).Value = 0
$ws.Cells.Item(26, 1).Value = "Bi1.5Ag0.5Ti2O7"
$ws.Cells.Item(26, 2).Value = 0
$ws.Cells.Item(27, 1).Value = "Bi1.5Mg0.5Ti2O7"
$ws.Cells.Item(27, 2).Value = 0
$ws.Cells.Item(28, 1).Value = "Bi1.5Ca0.5Ti2O7"
$ws.Cells.Item(28, 2).Value = 3.6
$ws.Cells.Item(29, 1).Value = "Bi1.5Zn0.5Ti2O7"
$ws.Cells.Item(29, 2).Value = 3.46
$ws.Cells.Item(30, 1).Value = "Bi1.5Cd0.5Ti2O7"
$ws.Cells.Item(30, 2).Value = 3.36
$ws.Cells.Item(31, 1).Value = "Bi1.5Cr0.5Ti2O7"
$ws.Cells.Item(31, 2).Value = 3.18
$ws.Cells.Item(32, 1).Value = "Bi1.5Mn0.5Ti2O7"
$ws.Cells.Item(32, 2).Value = 0
$ws.Cells.Item(33, 1).Value = "Bi1.5Fe0.5Ti2O7"
$ws.Cells.Item(33, 2).Value = 0
$ws.Cells.Item(34, 1).Value = "Bi1.5Co0.5Ti2O7"
$ws.Cells.Item(34, 2).Value = 0
$ws.Cells.Item(35, 1).Value = "Bi1.5Ni0.5Ti2O7"
$ws.Cells.Item(35, 2).Value = 0
$ws.Cells.Item(36, 1).Value = "Bi1.5Ce0.5Ti2O7"
$ws.Cells.Item(36, 2).Value = 0
$ws.Cells.Item(37, 1).Value = "Bi1.5Pr0.5Ti2O7"
$ws.Cells.Item(37, 2).Value = 0
$ws.Cells.Item(38, 1).Value = "Bi1.5Nd0.5Ti2O7"
$ws.Cells.Item(38, 2).Value = 0
$ws.Cells.Item(39, 1).Value = "Bi1.5Pm0.5Ti2O7"
$ws.Cells.Item(39, 2).Value = 0
$ws.Cells.Item(40, 1).Value = "Bi1.5Sm0.5Ti2O7"
$ws.Cells.Item(40, 2).Value = 0
$ws.Cells.Item(41, 1).Value = "Bi1.5Eu0.5Ti2O7"
$ws.Cells.Item(41, 2).Value = 0
$ws.Cells.Item(42, 1).Value = "Bi1.5Gd0.5Ti2O7"
$ws.Cells.Item(42, 2).Value = 0
$ws.Cells.Item(43, 1).Value = "Bi1.5Tb0.5Ti2O7"
$ws.Cells.Item(43, 2).Value = 0
$ws.Cells.Item(44, 1).Value = "Bi1.5Dy0.5Ti2O7"
$ws.Cells.Item(44, 2).Value = 0
$ws.Cells.Item(45, 1).Value = "Bi1.5Ho0.5Ti2O7"
$ws.Cells.Item(45, 2).Value = 0
$ws.Cells.Item(46, 1).Value = "Bi1.5Er0.5Ti2O7"
$ws.Cells.Item(46, 2).Value = 3.66
$ws.Cells.Item(47, 1).Value = "Bi1.5Tm0.5Ti2O7"
$ws.Cells.Item(47, 2).Value = 0
$ws.Cells.Item(48, 1).Value = "Bi1.5Yb0.5Ti2O7"
$ws.Cells.Item(48, 2).Value = 0
$ws.Cells.Item(49, 1).Value = "Bi1.5Lu0.5Ti2O7"
$ws.Cells.Item(49, 2).Value = 3.88
$ws.Cells.Item(50, 1).Value = "Bi1.5Eu0.5Ti2O7"
$ws.Cells.Item(50, 2).Value = 0
$ws.Cells.Item(51, 1).Value = "Bi1.75Eu0.25Ti2O7"
$ws.Cells.Item(51, 2).Value = 0
$ws.Cells.Item(52, 1).Value = "Bi1.875Eu0.125Ti2O7"
$ws.Cells.Item(52, 2).Value = 0
$ws.Cells.Item(53, 1).Value = "Bi1.5Ho0.5Ti2O7"
$ws.Cells.Item(53, 2).Value = 0
$ws.Cells.Item(54, 1).Value = "Bi1.75Ho0.25Ti2O7"
$ws.Cells.Item(54, 2).Value = 3.52
$ws.Cells.Item(55, 1).Value = "Bi1.875Ho0.125Ti2O7"
$ws.Cells.Item(55, 2).Value = 0
$ws.Cells.Item(56, 1).Value = "Bi1.5Yb0.5Ti2O7"
$ws.Cells.Item(56, 2).Value = 0
$ws.Cells.Item(57, 1).Value = "Bi1.75Yb0.25Ti2O7"
$ws.Cells.Item(57, 2).Value = 3.49
$ws.Cells.Item(58, 1).Value = "Bi1.875Yb0.125Ti2O7"
$ws.Cells.Item(58, 2).Value = 0
$ws.Cells.Item(59, 1).Value = "Bi1.6Li0.4Ti2O6.6"
$ws.Cells.Item(59, 2).Value = 3.84
$ws.Cells.Item(60, 1).Value = "Bi1.5Na0.125Ti2O6.3125"
$ws.Cells.Item(60, 2).Value = 3.66
$ws.Cells.Item(61, 1).Value = "Bi1.5Ga0.25Ti2O6.625"
$ws.Cells.Item(61, 2).Value = 3.66
$ws.Cells.Item(62, 1).Value = "Bi1.5In0.25Ti2O6.625"
$ws.Cells.Item(62, 2).Value = 3.53
$ws.Cells.Item(63, 1).Value = "Bi1.75In0.25Ti2O7"
$ws.Cells.Item(63, 2).Value = 3.49
$ws.Cells.Item(64, 1).Value = "Bi1.875In0.125Ti2O7"
$ws.Cells.Item(64, 2).Value = 3.48
$ws.Cells.Item(65, 1).Value = "Bi1.9375In0.0625Ti2O7"
$ws.Cells.Item(65, 2).Value = 3.47
$ws.Cells.Item(66, 1).Value = "Bi1.5Sc0.25Ti2O6.625"
$ws.Cells.Item(66, 2).Value = 3.63
$ws.Cells.Item(67, 1).Value = "Bi1.75Sc0.25Ti2O7"
$ws.Cells.Item(67, 2).Value = 3.57
$ws.Cells.Item(68, 1).Value = "Bi1.875Sc0.125Ti2O7"
$ws.Cells.Item(68, 2).Value = 3.48
$ws.Cells.Item(69, 1).Value = "Bi1.9375Sc0.0625Ti2O7"
$ws.Cells.Item(69, 2).Value = 3.36
$ws.Cells.Item(70, 1).Value = "Bi1.5Al0.25Ti2O6.625"
$ws.Cells.Item(70, 2).Value = 3.65
$ws.Cells.Item(71, 1).Value = "Bi1.75Al0.25Ti2O7"
$ws.Cells.Item(71, 2).Value = 3.55
$ws.Cells.Item(72, 1).Value = "Bi1.875Al0.125Ti2O7"
$ws.Cells.Item(72, 2).Value = 3.59
$ws.Cells.Item(73, 1).Value = "Bi1.9375Al0.0625Ti2O7"
$ws.Cells.Item(73, 2).Value = 3.46
$ws.Cells.Item(74, 1).Value = "Ta1.94Eu0.06O5"
$ws.Cells.Item(74, 2).Value = 2.72
$ws.Cells.Item(75, 1).Value = "BaTaO2N"
$ws.Cells.Item(75, 2).Value = 1.96
$ws.Cells.Item(76, 1).Value = "BaTa0.5Al0.5O2N"
$ws.Cells.Item(76, 2).Value = 2.31
$ws.Cells.Item(77, 1).Value = "BaTa0.5Mg0.5O2N"
$ws.Cells.Item(77, 2).Value = 2.8
$ws.Cells.Item(78, 1).Value = "BaTa0.5Al0.25Mg0.25O2N"
$ws.Cells.Item(78, 2).Value = 2.67
$ws.Cells.Item(79, 1).Value = "BaTa0.5Al0.125Mg0.375O2N"
$ws.Cells.Item(79, 2).Value = 2.68
$ws.Cells.Item(80, 1).Value = "BaTa0.5Al0.375Mg0.125O2N"
$ws.Cells.Item(80, 2).Value = 2.61
$ws.Cells.Item(81, 1).Value = "Al2O3"
$ws.Cells.Item(81, 2).Value = 6.75
$ws.Cells.Item(82, 1).Value = "TiO2"
$ws.Cells.Item(82, 2).Value = 0
$ws.Cells.Item(83, 1).Value = "Ta2O5"
$ws.Cells.Item(83, 2).Value = 3.48
$ws.Cells.Item(84, 1).Value = "Fe2O3"
$ws.Cells.Item(84, 2).Value = 0
$ws.Cells.Item(85, 1).Value = "NaLaTi2O6"
$ws.Cells.Item(85, 2).Value = 3.59
$ws.Cells.Item(86, 1).Value = "NaCl"
$ws.Cells.Item(86, 2).Value = 7.35
$ws.Cells.Item(87, 1).Value = "Bi1.5Mg1.0Nb1.5O7.0"
$ws.Cells.Item(87, 2).Value = 3.49
$ws.Cells.Item(88, 1).Value = "Bi1.5Mg0.9Na0.1Nb1.5O6.95"
$ws.Cells.Item(88, 2).Value = 3.59
$ws.Cells.Item(89, 1).Value = "Bi1.5Mg0.75Na0.25Nb1.5O6.875"
$ws.Cells.Item(89, 2).Value = 3.67
$ws.Cells.Item(90, 1).Value = "Bi1.5Mg0.65Na0.25Nb1.5O6.775"
$ws.Cells.Item(90, 2).Value = 3.57
$ws.Cells.Item(91, 1).Value = "Bi1.5Mg0.5Na0.4Nb1.5O6.7"
$ws.Cells.Item(91, 2).Value = 3.59
$ws.Cells.Item(92, 1).Value = "Bi1.5Mg0.65Li0.25Nb1.5O6.775"
$ws.Cells.Item(92, 2).Value = 3.41
$ws.Cells.Item(93, 1).Value = "Bi1.5Mg0.5Li0.4Nb1.5O6.7"
$ws.Cells.Item(93, 2).Value = 3.47
$ws.Cells.Item(94, 1).Value = "Bi1.0Eu0.5Li0.4Mg0.5Nb1.5O6.7"
$ws.Cells.Item(94, 2).Value = 0
$ws.Cells.Item(95, 1).Value = "Bi1.1Eu0.4Li0.4Mg0.5Nb1.5O6.7"
$ws.Cells.Item(95, 2).Value = 0
$ws.Cells.Item(96, 1).Value = "Bi1.2Eu0.3Li0.4Mg0.5Nb1.5O6.7"
$ws.Cells.Item(96, 2).Value = 0
$ws.Cells.Item(97, 1).Value = "Bi1.3Eu0.2Li0.4Mg0.5Nb1.5O6.7"
$ws.Cells.Item(97, 2).Value = 3.54
$ws.Cells.Item(98, 1).Value = "Bi1.4Eu0.1Li0.4Mg0.5Nb1.5O6.7"
$ws.Cells.Item(98, 2).Value = 3.38
$ws.Cells.Item(99, 1).Value = "Bi1.0Eu0.5Na0.4Mg0.5Nb1.5O6.7"
$ws.Cells.Item(99, 2).Value = 0
$ws.Cells.Item(100, 1).Value = "Bi1.1Eu0.4Na0.4Mg0.5Nb1.5O6.7"
$ws.Cells.Item(100, 2).Value = 0
$ws.Cells.Item(101, 1).Value = "Bi1.2Eu0.3Na0.4Mg0.5Nb1.5O6.7"
$ws.Cells.Item(101, 2).Value = 3.63
$ws.Cells.Item(102, 1).Value = "Bi1.3Eu0.2Na0.4Mg0.5Nb1.5O6.7"
$ws.Cells.Item(102, 2).Value = 3.58
$ws.Cells.Item(103, 1).Value = "Bi1.4Eu0.1Na0.4Mg0.5Nb1.5O6.7"
$ws.Cells.Item(103, 2).Value = 3.39
$ws.Cells.Item(104, 1).Value = "Bi1.5Mg0.9Nb1.5O6.9"
$ws.Cells.Item(104, 2).Value = 3.38
$ws.Cells.Item(105, 1).Value = "Bi1.5Mg0.5Al0.4Nb1.5O7.1"
$ws.Cells.Item(105, 2).Value = 3.26
$ws.Cells.Item(106, 1).Value = "Bi1.5Zn0.9Nb1.5O6.9"
$ws.Cells.Item(106, 2).Value = 3.25
$ws.Cells.Item(107, 1).Value = "Bi1.5Zn0.5Al0.4Nb1.5O7.1"
$ws.Cells.Item(107, 2).Value = 3.51
$ws.Cells.Item(108, 1).Value = "Bi1.3Li0.45Al0.45Nb1.5O6.6"
$ws.Cells.Item(108, 2).Value = 3.36
$ws.Cells.Item(109, 1).Value = "Bi1.3Li0.45Ga0.45Nb1.5O6.6"
$ws.Cells.Item(109, 2).Value = 3.46
$ws.Cells.Item(110, 1).Value = "Bi1.3Li0.45In0.45Nb1.5O6.6"
$ws.Cells.Item(110, 2).Value = 3.37
$ws.Cells.Item(111, 1).Value = "Bi1.5Mg0.5Ga0.4Nb1.5O7.1"
$ws.Cells.Item(111, 2).Value = 3.24
$ws.Cells.Item(112, 1).Value = "Bi1.5Mg0.5In0.4Nb1.5O7.1"
$ws.Cells.Item(112, 2).Value = 3.3
$ws.Cells.Item(113, 1).Value = "Bi1.5Zn0.5Ga0.4Nb1.5O7.1"
$ws.Cells.Item(113, 2).Value = 3.44
$ws.Cells.Item(114, 1).Value = "Bi1.5Zn0.5In0.4Nb1.5O7.1"
$ws.Cells.Item(114, 2).Value = 3.34
$ws.Cells.Item(115, 1).Value = "Bi1.3Na0.45Al0.45Nb1.5O6.6"
$ws.Cells.Item(115, 2).Value = 3.46
$ws.Cells.Item(116, 1).Value = "Bi1.3Na0.45Ga0.45Nb1.5O6.6"
$ws.Cells.Item(116, 2).Value = 3.46
$ws.Cells.Item(117, 1).Value = "Bi1.3Na0.45In0.45Nb1.5O6.6"
$ws.Cells.Item(117, 2).Value = 3.45
$ws.Cells.Item(118, 1).Value = "Bi1.56Cu0.19Mg0.58Nb1.56O7.01"
$ws.Cells.Item(118, 2).Value = 2.91
$ws.Cells.Item(119, 1).Value = "Bi1.56Cu0.39Mg0.39Nb1.56O7.02"
$ws.Cells.Item(119, 2).Value = 3.07
$ws.Cells.Item(120, 1).Value = "Bi1.46Li0.1Cu0.38Mg0.39Nb1.56O6.91"
$ws.Cells.Item(120, 2).Value = 2.71
$ws.Cells.Item(121, 1).Value = "Bi1.35Li0.21Cu0.38Mg0.39Nb1.56O6.8"
$ws.Cells.Item(121, 2).Value = 2.99
$ws.Cells.Item(122, 1).Value = "Bi1.5Li0.2Na0.2Mg0.5Nb1.5O6.7"
$ws.Cells.Item(122, 2).Value = 3.65
$ws.Cells.Item(123, 1).Value = "Bi1.3Li0.2Na0.2La0.2Mg0.5Nb1.5O6.7"
$ws.Cells.Item(123, 2).Value = 3.75
$ws.Cells.Item(124, 1).Value = "Bi1.3Li0.2Na0.2Eu0.2Mg0.5Nb1.5O6.7"
$ws.Cells.Item(124, 2).Value = 3.62
$ws.Cells.Item(125, 1).Value = "Bi1.1Li0.2Na0.2La0.2Eu0.2Mg0.5Nb1.5O6.7"
$ws.Cells.Item(125, 2).Value = 3.65
$ws.Cells.Item(126, 1).Value = "Bi0.9Li0.2Na0.2La0.3Eu0.3Mg0.5Nb1.5O6.7"
$ws.Cells.Item(126, 2).Value = 0
